$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rename the header strings: "_old" -> "_FV2304", "_new" -> "_FV2310"
#    (column K / "diff" is left untouched)
# ---------------------------------------------------------------------
$headers = @(
  "Segmentname_FV2304",
  "Segmentgruppe_FV2304",
  "Segment_FV2304",
  "Datenelement_FV2304",
  "Segment ID_FV2304",
  "Code_FV2304",
  "Qualifier_FV2304",
  "Beschreibung_FV2304",
  "Bedingungsausdruck_FV2304",
  "Bedingung_FV2304",
  "diff",
  "Segmentname_FV2310",
  "Segmentgruppe_FV2310",
  "Segment_FV2310",
  "Datenelement_FV2310",
  "Segment ID_FV2310",
  "Code_FV2310",
  "Qualifier_FV2310",
  "Beschreibung_FV2310",
  "Bedingungsausdruck_FV2310",
  "Bedingung_FV2310"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---------------------------------------------------------------------
# 2) Add a Table (ListObject) over A1:U58.
#    The existing header row already carries direct formatting
#    (bold / fill / border); if a table is built straight on top of it
#    Excel "bakes" that formatting into a headerRowDxf + dxfs entry,
#    which the target workbook does not have. To avoid mutating
#    xl/styles.xml or the header row's cell style, build the table on a
#    scratch, unformatted range first and then resize it onto the real
#    data - that path does not capture/convert any formatting.
# ---------------------------------------------------------------------
$scratch = $ws.Range("W1:X2")
$scratch.Cells.Item(1, 1).Value = "h1"
$scratch.Cells.Item(1, 2).Value = "h2"
$scratch.Cells.Item(2, 1).Value = "a"
$scratch.Cells.Item(2, 2).Value = "b"

$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $scratch, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = $null

$tbl.Resize($ws.Range("A1:U58"))

# Re-stamp the header row values so the table's column names pick up
# the renamed headers (Resize keeps the column names captured at
# creation time).
$hdrRange = $tbl.HeaderRowRange
for ($i = 0; $i -lt $headers.Length; $i++) {
  $hdrRange.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Drop the scratch cells used to seed the table.
$ws.Range("W1:X2").ClearContents()

# ---------------------------------------------------------------------
# 3) Freeze the header row (pane split after row 1).
# ---------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
